$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ptn"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.459193
$ws.Range("H2").Value = 4.377579000000001
$ws.Range("I2").Value = 0.01359248715138807
$ws.Range("J2").Value = 0.01359248715138807
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09880833333333333
$ws.Range("N2").Value = 0.296425
$ws.Range("O2").Value = 0.1308359775655526
$ws.Range("P2").Value = 0.1308359775655526
$ws.Range("Q2").Value = 0.1441804283416667
$ws.Range("R2").Value = 1.297623855075
$ws.Range("S2").Value = 0.001778386343999072
$ws.Range("T2").Value = 0.001778386343999071

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ptn"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.459193
$ws.Range("H3").Value = 4.377579000000001
$ws.Range("I3").Value = 0.01359248715138807
$ws.Range("J3").Value = 0.01359248715138807
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.05372733333333333
$ws.Range("N3").Value = 0.161182
$ws.Range("O3").Value = 0.07114246280162233
$ws.Range("P3").Value = 0.07114246280162234
$ws.Range("Q3").Value = 0.07839854870866667
$ws.Range("R3").Value = 0.7055869383780001
$ws.Range("S3").Value = 0.0009670030115491552
$ws.Range("T3").Value = 0.0009670030115491552

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ptn"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.459193
$ws.Range("H4").Value = 4.377579000000001
$ws.Range("I4").Value = 0.01359248715138807
$ws.Range("J4").Value = 0.01359248715138807
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.602672
$ws.Range("N4").Value = 1.808016
$ws.Range("O4").Value = 0.7980215596328251
$ws.Range("P4").Value = 0.7980215596328251
$ws.Range("Q4").Value = 0.8794147636960001
$ws.Range("R4").Value = 7.914732873264001
$ws.Range("S4").Value = 0.01084709779583984
$ws.Range("T4").Value = 0.01084709779583984

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ptn"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 103.0385286666667
$ws.Range("H5").Value = 309.115586
$ws.Range("I5").Value = 0.9598112634857745
$ws.Range("J5").Value = 0.9598112634857743
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09880833333333333
$ws.Range("N5").Value = 0.296425
$ws.Range("O5").Value = 0.1308359775655526
$ws.Range("P5").Value = 0.1308359775655526
$ws.Range("Q5").Value = 10.18106528667222
$ws.Range("R5").Value = 91.62958758005
$ws.Range("S5").Value = 0.1255778449365895
$ws.Range("T5").Value = 0.1255778449365895

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ptn"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 103.0385286666667
$ws.Range("H6").Value = 309.115586
$ws.Range("I6").Value = 0.9598112634857745
$ws.Range("J6").Value = 0.9598112634857743
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05372733333333333
$ws.Range("N6").Value = 0.161182
$ws.Range("O6").Value = 0.07114246280162233
$ws.Range("P6").Value = 0.07114246280162234
$ws.Range("Q6").Value = 5.535985375850221
$ws.Range("R6").Value = 49.823868382652
$ws.Range("S6").Value = 0.06828333710911484
$ws.Range("T6").Value = 0.06828333710911484

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ptn"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 103.0385286666667
$ws.Range("H7").Value = 309.115586
$ws.Range("I7").Value = 0.9598112634857745
$ws.Range("J7").Value = 0.9598112634857743
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.602672
$ws.Range("N7").Value = 1.808016
$ws.Range("O7").Value = 0.7980215596328251
$ws.Range("P7").Value = 0.7980215596328251
$ws.Range("Q7").Value = 62.09843614859733
$ws.Range("R7").Value = 558.885925337376
$ws.Range("S7").Value = 0.7659500814400702
$ws.Range("T7").Value = 0.76595008144007

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ptn"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.855184666666666
$ws.Range("H8").Value = 8.565553999999999
$ws.Range("I8").Value = 0.02659624936283746
$ws.Range("J8").Value = 0.02659624936283746
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.09880833333333333
$ws.Range("N8").Value = 0.296425
$ws.Range("O8").Value = 0.1308359775655526
$ws.Range("P8").Value = 0.1308359775655526
$ws.Range("Q8").Value = 0.2821160382722222
$ws.Range("R8").Value = 2.53904434445
$ws.Range("S8").Value = 0.003479746284964046
$ws.Range("T8").Value = 0.003479746284964045

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ptn"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.855184666666666
$ws.Range("H9").Value = 8.565553999999999
$ws.Range("I9").Value = 0.02659624936283746
$ws.Range("J9").Value = 0.02659624936283746
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.05372733333333333
$ws.Range("N9").Value = 0.161182
$ws.Range("O9").Value = 0.07114246280162233
$ws.Range("P9").Value = 0.07114246280162234
$ws.Range("Q9").Value = 0.1534014583142222
$ws.Range("R9").Value = 1.380613124828
$ws.Range("S9").Value = 0.001892122680958336
$ws.Range("T9").Value = 0.001892122680958336

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ptn"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.855184666666666
$ws.Range("H10").Value = 8.565553999999999
$ws.Range("I10").Value = 0.02659624936283746
$ws.Range("J10").Value = 0.02659624936283746
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.602672
$ws.Range("N10").Value = 1.808016
$ws.Range("O10").Value = 0.7980215596328251
$ws.Range("P10").Value = 0.7980215596328251
$ws.Range("Q10").Value = 1.720739853429333
$ws.Range("R10").Value = 15.486658680864
$ws.Range("S10").Value = 0.02122438039691508
$ws.Range("T10").Value = 0.02122438039691508
